$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Location")

# Colors used by the existing "Location" table (Numbers-exported banded
# rows): light-grey band (row 9/11/13 - "odd"), slightly darker band
# (row 10/12 - "even"), white thin gridlines, black 10pt Helvetica Neue text.
$bandOdd  = 13290186   # RGB CACACA, BGR-packed for OLE Color
$bandEven = 15132390   # RGB E6E6E6, BGR-packed for OLE Color
$gridLine = 16777215   # RGB FFFFFF
$textCol  = 0          # RGB 000000

# New test-case rows describing how to locate the 2nd/3rd rows of the tank
# list, and the description area, via XPath — added to support the new
# "select row on tank list" test cases.
$newRows = @(
    @{ Row = 9;  Key = "Home.List.All.Data2.Name";  Type = "XPath"; Value = "//*[@id='tanklist']/tbody/tr[2]/td[2]" },
    @{ Row = 10; Key = "Home.List.All.Data2.Stage"; Type = "XPath"; Value = "//*[@id='tanklist']/tbody/tr[2]/td[3]" },
    @{ Row = 11; Key = "Home.List.All.Data3.Name";  Type = "XPath"; Value = "//*[@id='tanklist']/tbody/tr[3]/td[2]" },
    @{ Row = 12; Key = "Home.List.All.Data3.Type";  Type = "XPath"; Value = "//*[@id='tanklist']/tbody/tr[3]/td[4]" },
    @{ Row = 13; Key = "Home.Desc.Name";            Type = "XPath"; Value = "//*[@id='imganddesc_div']/div/h4" }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    if (($r - 1) % 2 -eq 0) { $band = $bandOdd } else { $band = $bandEven }

    $ws.Rows.Item($r).RowHeight = 20.25

    # Key / Type / Value cells (columns A-C) carry the Text number format,
    # same as the rest of the table.
    $vals = @($entry.Key, $entry.Type, $entry.Value)
    for ($col = 1; $col -le 3; $col++) {
        $c = $ws.Cells.Item($r, $col)
        $c.Value = $vals[$col - 1]
        $c.NumberFormat = "@"
        $c.Font.Name = "Helvetica Neue"
        $c.Font.Size = 10
        $c.Font.Color = $textCol
        $c.Interior.Color = $band
        $c.Borders.Item(7).LineStyle = 1
        $c.Borders.Item(7).Color = $gridLine
        $c.Borders.Item(8).LineStyle = 1
        $c.Borders.Item(8).Color = $gridLine
        $c.Borders.Item(9).LineStyle = 1
        $c.Borders.Item(9).Color = $gridLine
        $c.Borders.Item(10).LineStyle = 1
        $c.Borders.Item(10).Color = $gridLine
        $c.VerticalAlignment = -4160
        $c.WrapText = $true
    }

    # Remaining columns (D-G) stay empty but keep the same banded
    # formatting as the rest of the table.
    for ($col = 4; $col -le 7; $col++) {
        $c = $ws.Cells.Item($r, $col)
        $c.Font.Name = "Helvetica Neue"
        $c.Font.Size = 10
        $c.Font.Color = $textCol
        $c.Interior.Color = $band
        $c.Borders.Item(7).LineStyle = 1
        $c.Borders.Item(7).Color = $gridLine
        $c.Borders.Item(8).LineStyle = 1
        $c.Borders.Item(8).Color = $gridLine
        $c.Borders.Item(9).LineStyle = 1
        $c.Borders.Item(9).Color = $gridLine
        $c.Borders.Item(10).LineStyle = 1
        $c.Borders.Item(10).Color = $gridLine
        $c.VerticalAlignment = -4160
        $c.WrapText = $true
    }
}

# Widen column A slightly to fit the new, longer keys (target raw OOXML
# width 25.5859 chars; Excel's ColumnWidth/raw-width pixel quantization
# means 24.85 is the closest achievable COM input).
$ws.Columns.Item(1).ColumnWidth = 24.85
